# Applies programmes_uqac (col J) and programmes_partenaires (col K) marker data
# to the "data" worksheet, plus the final cell selection/view state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 7764; $ws.Range("K2").Value = 1234
$ws.Range("J5").Value = "1307, 7764"; $ws.Range("K5").Value = 1234
$ws.Range("J9").Value = 7764; $ws.Range("K9").Value = 1234
$ws.Range("J11").Value = "3037,1709,3754,6908"
$ws.Range("J14").Value = "1307, 7764"; $ws.Range("K14").Value = 1234
$ws.Range("J16").Value = 1709.3753999999999
$ws.Range("J21").Value = "3037,1709,3754,6908"
$ws.Range("J26").Value = 1709.3753999999999
$ws.Range("J28").Value = 7764; $ws.Range("K28").Value = 1234
$ws.Range("J29").Value = "1307, 7764"; $ws.Range("K29").Value = 1234
$ws.Range("J33").Value = "3037,1709,3754,6908"
$ws.Range("J37").Value = 1709.3753999999999
$ws.Range("J40").Value = 7764; $ws.Range("K40").Value = 1234
$ws.Range("J44").Value = "1307, 7764"; $ws.Range("K44").Value = 1234
$ws.Range("J48").Value = "3037,1709,3754,6908"
$ws.Range("J54").Value = 1709.3753999999999
$ws.Range("J59").Value = "1307, 7764"; $ws.Range("K59").Value = 1234
$ws.Range("J62").Value = 7764; $ws.Range("K62").Value = 1234
$ws.Range("J63").Value = "3037,1709,3754,6908"
$ws.Range("J69").Value = 1709.3753999999999
$ws.Range("J74").Value = "1307, 7764"; $ws.Range("K74").Value = 1234
$ws.Range("J79").Value = 7764; $ws.Range("K79").Value = 1234
$ws.Range("J81").Value = "3037,1709,3754,6908"
$ws.Range("J85").Value = 1709.3753999999999
$ws.Range("J89").Value = "1307, 7764"
$ws.Range("J94").Value = "3037,1709,3754,6908"
$ws.Range("J97").Value = 7764; $ws.Range("K97").Value = 1234
$ws.Range("J100").Value = 1709.3753999999999
$ws.Range("J102").Value = "3037,1709,3754,6908"
$ws.Range("J107").Value = "1307, 7764"; $ws.Range("K107").Value = 1234
$ws.Range("J111").Value = "3037,1709,3754,6908"
$ws.Range("J115").Value = 7764; $ws.Range("K115").Value = 1234
$ws.Range("J117").Value = 1709.3753999999999
$ws.Range("J119").Value = "1307, 7764"; $ws.Range("K119").Value = 1234
$ws.Range("J126").Value = 1709.3753999999999
$ws.Range("J133").Value = 7764; $ws.Range("K133").Value = 1234
$ws.Range("J137").Value = "1307, 7764"; $ws.Range("K137").Value = 1234
$ws.Range("J140").Value = "3037,1709,3754,6908"
$ws.Range("J142").Value = 1709.3753999999999
$ws.Range("J145").Value = 7764; $ws.Range("K145").Value = 1234
$ws.Range("J153").Value = 1709.3753999999999
$ws.Range("J155").Value = "1307, 7764"; $ws.Range("K155").Value = 1234
$ws.Range("J158").Value = 1709.3753999999999
$ws.Range("J160").Value = 7764; $ws.Range("K160").Value = 1234
$ws.Range("J163").Value = "3037,1709,3754,6908"
$ws.Range("J165").Value = 1709.3753999999999
$ws.Range("J170").Value = "1307, 7764"; $ws.Range("K170").Value = 1234
$ws.Range("J176").Value = 1709.3753999999999
$ws.Range("J178").Value = 7764; $ws.Range("K178").Value = 1234
$ws.Range("J183").Value = "3037,1709,3754,6908"
$ws.Range("J185").Value = "1307, 7764"; $ws.Range("K185").Value = 1234
$ws.Range("J187").Value = 1709.3753999999999
$ws.Range("J196").Value = 7764; $ws.Range("K196").Value = 1234
$ws.Range("J198").Value = 1709.3753999999999
$ws.Range("J203").Value = "1307, 7764"; $ws.Range("K203").Value = 1234
$ws.Range("J207").Value = "3037,1709,3754,6908"
$ws.Range("J213").Value = 1709.3753999999999
$ws.Range("J217").Value = 7764; $ws.Range("K217").Value = 1234
$ws.Range("J221").Value = "1307, 7764"; $ws.Range("K221").Value = 1234
$ws.Range("J227").Value = "3037,1709,3754,6908"
$ws.Range("J229").Value = 1709.3753999999999
$ws.Range("J233").Value = 7764; $ws.Range("K233").Value = 1234
$ws.Range("J234").Value = "1307, 7764"; $ws.Range("K234").Value = 1234

# Restore the view/selection state recorded at save time
$win = $excel.ActiveWindow
$win.ScrollRow = 202
$win.ScrollColumn = 5
$ws.Range("K237").Select() | Out-Null
